$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Handout master: update the "Updated:" date placeholder text
#    11/19/20 -> 1/12/21
# ---------------------------------------------------------------------------
$hm = $p.HandoutMaster
$hf = $hm.HeadersFooters
$dt = $hf.DateAndTime
$dt.Text = "1/12/21"

# ---------------------------------------------------------------------------
# Locate slide 3 (the high-level flow diagram with the numbered callouts)
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(3)

# ---------------------------------------------------------------------------
# 2) Fix typo in the "Consortium aggres on CC Def" callout
#    aggres -> agrees
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        $tr = $shape.TextFrame.TextRange
        $full = $tr.Text
        $oldText = "Consortium aggres on CC Def"
        $idx = $full.IndexOf($oldText)
        if ($idx -ge 0) {
            $sub = $tr.Characters($idx + 1, $oldText.Length)
            $sub.Text = "Consortium agrees on CC Def"
        }
    }
}

# ---------------------------------------------------------------------------
# 3) Renumber the "Submit registerEnclave tx for ordering" callout
#    "12. Submit registerEnclave tx for ordering"
#      -> "13. " (new run) + "Submit registerEnclave tx for ordering"
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        $tr = $shape.TextFrame.TextRange
        $full = $tr.Text
        $oldPrefix = "12. Submit registerEnclave tx for ordering"
        $idx = $full.IndexOf($oldPrefix)
        if ($idx -ge 0) {
            $sub = $tr.Characters($idx + 1, 4)
            $sub.Text = "13. "
        }
    }
}
